$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.095.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.812.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "700.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.811.23"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("E11").Value = "  +4.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("D12").Style = "Normal"
$ws.Range("E13").Value = "  +8.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.457.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.811.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.083.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.51%  "
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +16.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "480.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.713"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.48%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.92%  "
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.966.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +15.35%  "
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("E33").Value = "  +5.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.187"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "29.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.86%  "
$ws.Range("E36").Value = "  +4.70%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  +2.90%  "
$ws.Range("E39").Value = "  +5.65%  "
$ws.Range("E40").Value = "  +2.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.979"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.21%  "
$ws.Range("E43").Value = "  +21.09%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "49.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "44.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.84%  "
$ws.Range("E49").Value = "  +1.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "410.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.06%  "
